$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "'" + "61.398.10"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value2 = "'" + "  +1.26%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value2 = "'" + "2.382.11"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value2 = "'" + "  +1.51%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value2 = "'" + "  -0.14%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value2 = "'" + "551.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value2 = "'" + "  +1.78%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value2 = "'" + "140.18"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value2 = "'" + "  +1.33%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D8").Value2 = "'" + "0.526"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value2 = "'" + "  +1.71%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value2 = "'" + "2.383.15"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value2 = "'" + "  +1.53%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value2 = "'" + "  +4.27%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value2 = "'" + "  +1.62%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value2 = "'" + "5.35"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value2 = "'" + "  +2.64%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value2 = "'" + "0.350"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value2 = "'" + "  +3.34%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value2 = "'" + "25.41"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value2 = "'" + "  +2.82%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value2 = "'" + "  +4.46%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value2 = "'" + "61.196.47"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value2 = "'" + "  +0.57%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value2 = "'" + "2.383.83"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value2 = "'" + "  +1.43%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value2 = "'" + "10.99"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value2 = "'" + "  +4.29%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value2 = "'" + "322.08"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value2 = "'" + "  +2.87%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value2 = "'" + "4.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value2 = "'" + "  +2.19%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value2 = "'" + "6.77"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value2 = "'" + "  +3.24%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value2 = "'" + "  +0.27%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value2 = "'" + "64.38"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value2 = "'" + "  +2.03%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value2 = "'" + "1.71"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value2 = "'" + "  -6.92%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value2 = "'" + "8.77"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value2 = "'" + "  +9.47%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value2 = "'" + "8.22"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value2 = "'" + "  +4.08%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value2 = "'" + "516.25"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value2 = "'" + "  +1.36%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value2 = "'" + "0.0$([char]8323)0903"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value2 = "'" + "  +0.98%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value2 = "'" + "  +5.20%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value2 = "'" + "1.39"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value2 = "'" + "  +1.57%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value2 = "'" + "  +1.87%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value2 = "'" + "1.55"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value2 = "'" + "  +1.49%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value2 = "'" + "  -0.07%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("B34").Value2 = "NEARProtocol"
$ws.Range("C34").Value2 = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").Value2 = "'" + "4.72"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value2 = "'" + "  +3.59%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("B35").Value2 = "RenderToken"
$ws.Range("C35").Value2 = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D35").Value2 = "'" + "5.52"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value2 = "'" + "  +5.14%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("B36").Value2 = "Stacks"
$ws.Range("C36").Value2 = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D36").Value2 = "'" + "1.91"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value2 = "'" + "  +6.69%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value2 = "'" + "  +2.75%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value2 = "'" + "18.55"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value2 = "'" + "  +2.45%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value2 = "'" + "147.08"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value2 = "'" + "  +5.68%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value2 = "'" + "  -0.10%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value2 = "'" + "41.26"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value2 = "'" + "  +2.91%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value2 = "'" + "149.99"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value2 = "'" + "  +8.99%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value2 = "'" + "2.17"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value2 = "'" + "  +3.69%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value2 = "'" + "3.61"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value2 = "'" + "  +2.10%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value2 = "'" + "  +3.08%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value2 = "'" + "19.62"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value2 = "'" + "  +1.35%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value2 = "'" + "0.580"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value2 = "'" + "  +2.71%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value2 = "'" + "0.0908"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value2 = "'" + "  +1.67%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value2 = "'" + "0.0225"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value2 = "'" + "  +1.92%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value2 = "'" + "11.41"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value2 = "'" + "  +0.60%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value2 = "'" + "16.81"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value2 = "'" + "  +1.31%  "
$ws.Range("E51").Style = "Normal"
